$d = $word.ActiveDocument

$replacements = @(
    @("623×3=1869", "130×9=1170"),
    @("127×3=381",  "379×7=2653"),
    @("924×5=4620", "950×5=4750"),
    @("554×4=2216", "364×6=2184"),
    @("112×8=896",  "622×2=1244"),
    @("627×7=4389", "409×3=1227"),
    @("194×7=1358", "614×8=4912"),
    @("357×5=1785", "951×3=2853"),
    @("388×3=1164", "364×2=728"),
    @("329×5=1645", "965×5=4825"),
    @("777×8=6216", "828×7=5796"),
    @("801×7=5607", "457×7=3199"),
    @("425×5=2125", "587×9=5283"),
    @("643×7=4501", "570×4=2280"),
    @("667×6=4002", "978×6=5868"),
    @("821×2=1642", "198×6=1188"),
    @("734×7=5138", "536×8=4288"),
    @("521×3=1563", "989×8=7912"),
    @("438×4=1752", "253×7=1771"),
    @("288×9=2592", "814×2=1628"),
    @("666×7=4662", "138×7=966"),
    @("919×5=4595", "506×2=1012"),
    @("766×5=3830", "658×5=3290"),
    @("754×5=3770", "267×7=1869"),
    @("361×3=1083", "960×7=6720")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
